# Refresh currentAveragePrice / LevePrice / LeveProfit columns (H:N) across all
# crafting-job leve tables, per the latest market snapshot from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3425
$ws.Range("I18").Value = 3425
$ws.Range("K18").Value = 3425
$ws.Range("M18").Value = -3141

$ws.Range("H58").Value = 1671.2858
$ws.Range("I58").Value = 36.4
$ws.Range("K58").Value = 109.2
$ws.Range("M58").Value = 40.80000000000001

$ws.Range("H112").Value = 2216.7856
$ws.Range("I112").Value = 945.4
$ws.Range("K112").Value = 2836.2
$ws.Range("M112").Value = -1728.2

$ws.Range("H116").Value = 3800.9048
$ws.Range("I116").Value = 3116.25
$ws.Range("J116").Value = 4222.231
$ws.Range("K116").Value = 3116.25
$ws.Range("L116").Value = 4222.231
$ws.Range("M116").Value = 325.75
$ws.Range("N116").Value = -11106.231

$ws.Range("H132").Value = 1681.6666
$ws.Range("I132").Value = 1715.9
$ws.Range("K132").Value = 5147.700000000001
$ws.Range("M132").Value = -2617.700000000001

$ws.Range("H141").Value = 4134.2
$ws.Range("I141").Value = 1993.5
$ws.Range("K141").Value = 5980.5
$ws.Range("M141").Value = -800.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1980
$ws.Range("I5").Value = 2425
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 2425
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -2313
$ws.Range("N5").Value = -424

$ws.Range("H61").Value = 10739.8
$ws.Range("I61").Value = 999.5
$ws.Range("K61").Value = 999.5
$ws.Range("M61").Value = -787.5

$ws.Range("H74").Value = 922.4545000000001
$ws.Range("I74").Value = 899.6667
$ws.Range("J74").Value = 949.8
$ws.Range("K74").Value = 899.6667
$ws.Range("L74").Value = 949.8
$ws.Range("M74").Value = -25.66669999999999
$ws.Range("N74").Value = -2697.8

$ws.Range("H77").Value = 922.4545000000001
$ws.Range("I77").Value = 899.6667
$ws.Range("J77").Value = 949.8
$ws.Range("K77").Value = 4498.3335
$ws.Range("L77").Value = 4749
$ws.Range("M77").Value = -130.3334999999997
$ws.Range("N77").Value = -13485

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""

$ws.Range("H136").Value = 10739.8
$ws.Range("I136").Value = 999.5
$ws.Range("K136").Value = 2998.5
$ws.Range("M136").Value = -448.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1980
$ws.Range("I4").Value = 2425
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 2425
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -2310
$ws.Range("N4").Value = -430

$ws.Range("H22").Value = 2067
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 201
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 201
$ws.Range("M22").Value = -2827
$ws.Range("N22").Value = -547

$ws.Range("H99").Value = 770
$ws.Range("I99").Value = 770
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 770
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 728
$ws.Range("N99").Value = ""

$ws.Range("H107").Value = 3364.5789
$ws.Range("I107").Value = 3032.6428
$ws.Range("K107").Value = 3032.6428
$ws.Range("M107").Value = -1112.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 908.3077
$ws.Range("I31").Value = 1115.5
$ws.Range("K31").Value = 1115.5
$ws.Range("M31").Value = -820.5

$ws.Range("H33").Value = 3594.5264
$ws.Range("I33").Value = 2137
$ws.Range("J33").Value = 4654.5454
$ws.Range("K33").Value = 2137
$ws.Range("L33").Value = 4654.5454
$ws.Range("M33").Value = -1758
$ws.Range("N33").Value = -5412.5454

$ws.Range("H34").Value = 908.3077
$ws.Range("I34").Value = 1115.5
$ws.Range("K34").Value = 1115.5
$ws.Range("M34").Value = -913.5

$ws.Range("H58").Value = 3358.75
$ws.Range("I58").Value = 1900.4
$ws.Range("K58").Value = 1900.4
$ws.Range("M58").Value = -1697.4

$ws.Range("H94").Value = 2139.4
$ws.Range("J94").Value = 969
$ws.Range("L94").Value = 969
$ws.Range("N94").Value = -1871

$ws.Range("H122").Value = 676.625
$ws.Range("I122").Value = 529.1111
$ws.Range("K122").Value = 1587.3333
$ws.Range("M122").Value = 862.6667000000002

$ws.Range("H132").Value = 5893.1333
$ws.Range("I132").Value = 5893.1333
$ws.Range("K132").Value = 17679.3999
$ws.Range("M132").Value = -15149.3999

$ws.Range("H136").Value = 3358.75
$ws.Range("I136").Value = 1900.4
$ws.Range("K136").Value = 5701.200000000001
$ws.Range("M136").Value = -3151.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 10607.031
$ws.Range("I56").Value = 10607.031
$ws.Range("K56").Value = 10607.031
$ws.Range("M56").Value = -10077.031

$ws.Range("H68").Value = 3933.9688
$ws.Range("J68").Value = 4022.1614
$ws.Range("L68").Value = 12066.4842
$ws.Range("N68").Value = -13688.4842

$ws.Range("H71").Value = 3933.9688
$ws.Range("J71").Value = 4022.1614
$ws.Range("L71").Value = 36199.4526
$ws.Range("N71").Value = -44311.4526

$ws.Range("H80").Value = 2151
$ws.Range("J80").Value = 2300
$ws.Range("L80").Value = 6900
$ws.Range("N80").Value = -8772

$ws.Range("H83").Value = 2151
$ws.Range("J83").Value = 2300
$ws.Range("L83").Value = 20700
$ws.Range("N83").Value = -30060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 58333.332
$ws.Range("I46").Value = 58000
$ws.Range("K46").Value = 58000
$ws.Range("M46").Value = -57844

$ws.Range("H47").Value = 35000
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36136

$ws.Range("H102").Value = 1814.579
$ws.Range("I102").Value = 1792
$ws.Range("J102").Value = 1899.25
$ws.Range("K102").Value = 1792
$ws.Range("L102").Value = 1899.25
$ws.Range("M102").Value = -170
$ws.Range("N102").Value = -5143.25

$ws.Range("H113").Value = 1985
$ws.Range("I113").Value = 1453.1666
$ws.Range("K113").Value = 1453.1666
$ws.Range("M113").Value = 716.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7057.8
$ws.Range("I7").Value = 5820.5
$ws.Range("K7").Value = 5820.5
$ws.Range("M7").Value = -5708.5

$ws.Range("H40").Value = 3008.0715
$ws.Range("I40").Value = 2607.1667
$ws.Range("K40").Value = 2607.1667
$ws.Range("M40").Value = -2471.1667

$ws.Range("H55").Value = 1320.5714
$ws.Range("I55").Value = 1232.5714
$ws.Range("K55").Value = 1232.5714
$ws.Range("M55").Value = -1059.5714

$ws.Range("H126").Value = 7057.8
$ws.Range("I126").Value = 5820.5
$ws.Range("K126").Value = 17461.5
$ws.Range("M126").Value = -14991.5

$ws.Range("H136").Value = 1500.0834
$ws.Range("I136").Value = 1525.875
$ws.Range("K136").Value = 4577.625
$ws.Range("M136").Value = -2027.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3823
$ws.Range("J96").Value = 3299
$ws.Range("L96").Value = 3299
$ws.Range("N96").Value = -6045

$ws.Range("H107").Value = 1084.421
$ws.Range("I107").Value = 957.8461
$ws.Range("J107").Value = 1358.6666
$ws.Range("K107").Value = 2873.5383
$ws.Range("L107").Value = 4075.9998
$ws.Range("M107").Value = -953.5383000000002
$ws.Range("N107").Value = -7915.9998

$ws.Range("H113").Value = 346.2
$ws.Range("I113").Value = 309.41666
$ws.Range("J113").Value = 493.33334
$ws.Range("K113").Value = 928.2499799999999
$ws.Range("L113").Value = 1480.00002
$ws.Range("M113").Value = 1241.75002
$ws.Range("N113").Value = -5820.000019999999

$ws.Range("H122").Value = 2556.3157
$ws.Range("I122").Value = 1953.8182
$ws.Range("J122").Value = 3384.75
$ws.Range("K122").Value = 5861.4546
$ws.Range("L122").Value = 10154.25
$ws.Range("M122").Value = -3411.4546
$ws.Range("N122").Value = -15054.25

$ws.Range("H136").Value = 2416.8147
$ws.Range("I136").Value = 2317.5
$ws.Range("K136").Value = 6952.5
$ws.Range("M136").Value = -4402.5
